# Diary.docx edit: "Grab all image data in one getImageData call rather
# than one for each pixel" — appends a new "Week 11" diary entry (a
# Heading 1 paragraph followed by a bulleted list) after the last
# paragraph of the document, using a brand new bulleted-list definition.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. Create a brand new bulleted-list definition (this is what Word
#    does internally whenever "Bullets" is applied to a paragraph that
#    isn't already part of a list) - this mints a new numId (10) in
#    numbering.xml that the new bullet paragraphs below will use.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$scratchPara = $d.Paragraphs($d.Paragraphs.Count)
$bulletTemplate = $word.ListGalleries(1).ListTemplates(1)
$scratchPara.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate)

# ---------------------------------------------------------------------
# 2. Build the replacement content: the existing last paragraph ("It is
#    therefore ...", stripped of the _GoBack bookmark it currently
#    carries) followed by the new "Week 11" heading, the new bulleted
#    paragraphs (all on the freshly minted numId), and the _GoBack
#    bookmark relocated onto the very last new paragraph - exactly
#    where Word leaves it after you type new content at the end of the
#    document.
# ---------------------------------------------------------------------
$weekHeading = "Week 11 (21/3/2016 " + [char]0x2013 + " 27/3/2016)"

$bulletPoints = @(
    "Began using Google Chrome" + [char]0x2019 + "s profiler to determine which parts of my program are taking the longest",
    "I started by trying to optimise the initial uploading of an image and extracting colours",
    "I soon discovered that it was the getImageData that was taking up the most processor time: 24,000 to 25,000 ms for a 350x235 image",
    "I tried sampling a quarter of all pixels rather than every one: while this brought the computation down to about 6,500ms it wasn" + [char]0x2019 + "t ideal and resulted in blocky images",
    "After a little Googling I realised that it was the getImageData call on the DOM object that was taking the longest, so I decided to get all the pixels in one getImageData call and iterate over the resulting array",
    "This brought the computation down to a blinding 66.1ms",
    "Hopefully I can do a similar operation for the putImageData method call",
    "I will also look into what processes are the most expensive when adjusting an image" + [char]0x2019 + "s colour"
)

$origAttrs = 'w14:paraId="6BB130EE" w14:textId="62AA6892" w:rsidR="00AF56C7" w:rsidRPr="00650226" w:rsidRDefault="00AF56C7" w:rsidP="00650226"'
$origNs = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$xml = "<w:p $wNs $origNs $origAttrs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""9""/></w:numPr></w:pPr><w:r><w:t>It is therefore even more important that I now look at how I can make this more efficient</w:t></w:r></w:p>"

$xml += "<w:p $wNs><w:pPr><w:pStyle w:val=""Heading1""/></w:pPr><w:r><w:t>$weekHeading</w:t></w:r></w:p>"

for ($i = 0; $i -lt $bulletPoints.Count; $i++) {
    $bookmark = ""
    if ($i -eq $bulletPoints.Count - 1) {
        $bookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }
    $xml += "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""10""/></w:numPr></w:pPr><w:r><w:t>$($bulletPoints[$i])</w:t></w:r>$bookmark</w:p>"
}

# ---------------------------------------------------------------------
# 3. Replace the span from the original last paragraph through the
#    scratch paragraph (the one that minted numId 10) with the XML
#    built above.
# ---------------------------------------------------------------------
$origPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$scratchParaEnd = $d.Paragraphs($d.Paragraphs.Count)
$target = $d.Range($origPara.Range.Start, $scratchParaEnd.Range.End)
$target.InsertXML($xml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
